{"js": "// Add a new \"Possible Destin Meetings\" paragraph at the end of the document,\n// right after the \"Contact ... for more information.\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Full paragraph text (heading + body sentence) inserted as one run first;\n// we will then re-format the \"Possible Destin\" / \" Meetings\" portions.\nconst headingPart1 = \"Possible Destin\";\nconst headingPart2 = \" Meetings\";\nconst bodyText =\n  \"  The Destin Library has nice facilities for small groups of  \" +\n  \"participants in a computer class.  Anyone interested in attending such \" +\n  \"a class for Perl programming please contact Tom.\";\n\nconst newParagraph = lastParagraph.insertParagraph(\n  headingPart1 + headingPart2 + bodyText,\n  \"After\"\n);\nawait context.sync();\n\n// Bold + single-underline the heading, split into two runs so the produced\n// OOXML matches \"Possible Destin\" and \" Meetings\" as distinct runs.\nconst headingRange1 = newParagraph.search(headingPart1, { matchCase: true });\nheadingRange1.load(\"items\");\nawait context.sync();\nheadingRange1.items[0].font.set({ bold: true, underline: \"Single\" });\nawait context.sync();\n\nconst headingRange2 = newParagraph.search(headingPart2, { matchCase: true });\nheadingRange2.load(\"items\");\nawait context.sync();\nheadingRange2.items[0].font.set({ bold: true, underline: \"Single\" });\nawait context.sync();\n", "ps1": "# Add a new \"Possible Destin Meetings\" paragraph at the end of the document,\n# right after the \"Contact ... for more information.\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the last paragraph in the document and insert a new, empty\n# paragraph immediately after it (inherits the sz/szCs=28 paragraph-mark\n# formatting from the preceding paragraph, same as Word does on Enter).\n$lastParaIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastParaIndex)\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n\n$headingPart1 = \"Possible Destin\"\n$headingPart2 = \" Meetings\"\n$bodyText = \"  The Destin Library has nice facilities for small groups of  participants in a computer class.  Anyone interested in attending such a class for Perl programming please contact Tom.\"\n\n# Fill in the whole paragraph text as plain text first.\n$newPara.Range.Text = $headingPart1 + $headingPart2 + $bodyText\n\n# Bold + single-underline \"Possible Destin\" (its own run).\n$find1 = $newPara.Range.Duplicate\n$find1.Find.Text = $headingPart1\n[void]$find1.Find.Execute()\n$find1.Bold = 1\n$find1.Underline = 1\n\n# Bold + single-underline \" Meetings\" (separate run from the above).\n$find2 = $newPara.Range.Duplicate\n$find2.Find.Text = $headingPart2\n[void]$find2.Find.Execute()\n$find2.Bold = 1\n$find2.Underline = 1\n"}
